$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.308.64"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.874.29"
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'0.7100"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "'241.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07793"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "'0.3106"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'25.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").Value = "'0.08434"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.860.56"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'5.237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "'0.7111"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "'91.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "29.317.03"
$ws.Range("D17").Value = "'0.000008311"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.51%  "
$ws.Range("D18").Value = "'6.071"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "'240.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "'13.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "2.111.86"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'0.1593"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").Value = "'162.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "'9.018"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'18.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'1.506"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "'4.402"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'1.300"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").Value = "'0.05364"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.84%  "
$ws.Range("D34").Value = "'1.944"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").Value = "'1.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "'0.7500"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("D37").Value = "'2.693"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "'0.01880"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "1.225.24"
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("D40").Value = "'2.730"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "'6.480"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").Value = "'0.8903"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'72.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").Value = "'108.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "2.011.46"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000123"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5201"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "'1.794"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'9.436"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'0.4323"
$ws.Range("D51").Style = "Normal"
